# This presentation's single slide contains three groups of shapes that were
# each accidentally duplicated (e.g. "Index Data" / "Market Internals Data",
# "Technical Indicators", and "Sentiment Data" / "Economic Data" sections).
# The fix removes the redundant duplicate copies and renumbers/renames the
# surviving shapes so the naming stays sequential, matching the v3 template
# cleanup ("feat: v3 templates with native shapes, editable text, 32 element
# PNGs from InDesign").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Delete the duplicate shapes (by their original auto-generated name) ---
$namesToDelete = @(
    "TextBox 10",
    "TextBox 11",
    "Rounded Rectangle 12",
    "TextBox 13",
    "TextBox 14",
    "Rounded Rectangle 15",
    "TextBox 19",
    "TextBox 20",
    "Rounded Rectangle 21",
    "TextBox 28",
    "TextBox 29",
    "Rounded Rectangle 30",
    "TextBox 31",
    "TextBox 32",
    "Rounded Rectangle 33"
)

foreach ($name in $namesToDelete) {
    $s.Shapes.Item($name).Delete()
}

# --- 2. Rename the surviving duplicate copies so naming stays sequential ---
$renames = @{
    "TextBox 16"            = "TextBox 10";
    "TextBox 17"             = "TextBox 11";
    "Rounded Rectangle 18"   = "Rounded Rectangle 12";
    "TextBox 22"             = "TextBox 13";
    "TextBox 23"             = "TextBox 14";
    "Rounded Rectangle 24"   = "Rounded Rectangle 15";
    "TextBox 25"             = "TextBox 16";
    "TextBox 26"             = "TextBox 17";
    "Rounded Rectangle 27"   = "Rounded Rectangle 18";
    "Rounded Rectangle 34"   = "Rounded Rectangle 19";
    "Rounded Rectangle 35"   = "Rounded Rectangle 20";
    "Rounded Rectangle 36"   = "Rounded Rectangle 21";
    "Rounded Rectangle 37"   = "Rounded Rectangle 22"
}

foreach ($oldName in $renames.Keys) {
    $s.Shapes.Item($oldName).Name = $renames[$oldName]
}
